$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8
$ws.Range("B8").Value = 0.7838
$ws.Range("D8").Value = 0.9941
$ws.Range("F8").Value = 0.997
$ws.Range("H8").Value = 0.3385
$ws.Range("J8").Value = 0.0315
$ws.Range("L8").Value = 0.8972
$ws.Range("N8").Value = 0.030411
$ws.Range("P8").Value = 0.2367

# Row 14
$ws.Range("B14").Value = 0.8132
$ws.Range("D14").Value = 1.0622
$ws.Range("F14").Value = 1.0306
$ws.Range("H14").Value = 0.2109
$ws.Range("J14").Value = 0.0291
$ws.Range("L14").Value = 0.9043
$ws.Range("N14").Value = 0.028193
$ws.Range("P14").Value = 0.2641

# Row 16
$ws.Range("B16").Value = 0.7847
$ws.Range("D16").Value = 0.9951
$ws.Range("F16").Value = 0.9976
$ws.Range("H16").Value = 0.2325
$ws.Range("J16").Value = 0.0224
$ws.Range("L16").Value = 0.453
$ws.Range("N16").Value = 0.021334
$ws.Range("P16").Value = 0.2385

# Row 21
$ws.Range("B21").Value = 0.7715
$ws.Range("D21").Value = 0.9732
$ws.Range("F21").Value = 0.9865
$ws.Range("H21").Value = 0.3532
$ws.Range("J21").Value = 0.032
$ws.Range("L21").Value = 0.8998
$ws.Range("N21").Value = 0.030891
$ws.Range("P21").Value = 0.2371

# Row 23
$ws.Range("B23").Value = 0.7715
$ws.Range("D23").Value = 0.9738
$ws.Range("F23").Value = 0.9868
$ws.Range("H23").Value = 0.3566
$ws.Range("J23").Value = 0.0318
$ws.Range("L23").Value = 0.899
$ws.Range("N23").Value = 0.030755
$ws.Range("P23").Value = 0.3068

# Row 40
$ws.Range("B40").Value = 0.7825
$ws.Range("D40").Value = 0.9909
$ws.Range("F40").Value = 0.9955
$ws.Range("H40").Value = 0.3422
$ws.Range("J40").Value = 0.0316
$ws.Range("L40").Value = 0.8971
$ws.Range("N40").Value = 0.030483
$ws.Range("P40").Value = 0.3449

# Update the active selection to match the saved state in the diff
$ws.Range("N41").Select()
